$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 29500
$ws.Range("I16").Value = 29500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 29500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -29270
$ws.Range("N16").ClearContents()

$ws.Range("H17").Value = 1962.1875
$ws.Range("I17").Value = 986
$ws.Range("K17").Value = 2958
$ws.Range("M17").Value = -2790

$ws.Range("H40").Value = 4027
$ws.Range("I40").Value = 3250.25
$ws.Range("J40").Value = 4415.375
$ws.Range("K40").Value = 3250.25
$ws.Range("L40").Value = 4415.375
$ws.Range("M40").Value = -3075.25
$ws.Range("N40").Value = -4765.375

$ws.Range("H112").Value = 1328
$ws.Range("J112").Value = 1425.6666
$ws.Range("L112").Value = 4276.9998
$ws.Range("N112").Value = -6492.9998

$ws.Range("H132").Value = 2003145.5
$ws.Range("I132").Value = 1133.4681
$ws.Range("J132").Value = 33368000
$ws.Range("K132").Value = 3400.4043
$ws.Range("L132").Value = 100104000
$ws.Range("M132").Value = -870.4043000000001
$ws.Range("N132").Value = -100109060

$ws.Range("H135").Value = 934.5
$ws.Range("I135").Value = 934.5
$ws.Range("K135").Value = 8410.5
$ws.Range("M135").Value = -5875.5

$ws.Range("H139").Value = 109619.71
$ws.Range("J139").Value = 109619.71
$ws.Range("L139").Value = 109619.71
$ws.Range("N139").Value = -119899.71

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2934.6667
$ws.Range("I74").Value = 2568.8462
$ws.Range("J74").Value = 3529.125
$ws.Range("K74").Value = 2568.8462
$ws.Range("L74").Value = 3529.125
$ws.Range("M74").Value = -1694.8462
$ws.Range("N74").Value = -5277.125

$ws.Range("H77").Value = 2934.6667
$ws.Range("I77").Value = 2568.8462
$ws.Range("J77").Value = 3529.125
$ws.Range("K77").Value = 12844.231
$ws.Range("L77").Value = 17645.625
$ws.Range("M77").Value = -8476.231
$ws.Range("N77").Value = -26381.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6031
$ws.Range("I31").Value = 2513.2222
$ws.Range("J31").Value = 7893.353
$ws.Range("K31").Value = 2513.2222
$ws.Range("L31").Value = 7893.353
$ws.Range("M31").Value = -2218.2222
$ws.Range("N31").Value = -8483.352999999999

$ws.Range("H34").Value = 6031
$ws.Range("I34").Value = 2513.2222
$ws.Range("J34").Value = 7893.353
$ws.Range("K34").Value = 2513.2222
$ws.Range("L34").Value = 7893.353
$ws.Range("M34").Value = -2311.2222
$ws.Range("N34").Value = -8297.352999999999

$ws.Range("H58").Value = 2129
$ws.Range("I58").Value = 2133.375
$ws.Range("K58").Value = 2133.375
$ws.Range("M58").Value = -1930.375

$ws.Range("H136").Value = 2129
$ws.Range("I136").Value = 2133.375
$ws.Range("K136").Value = 6400.125
$ws.Range("M136").Value = -3850.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 601.13336
$ws.Range("I5").Value = 601.13336
$ws.Range("K5").Value = 1803.40008
$ws.Range("M5").Value = -1691.40008

$ws.Range("H113").Value = 2430.3333
$ws.Range("I113").Value = 3150
$ws.Range("J113").Value = 991
$ws.Range("K113").Value = 9450
$ws.Range("L113").Value = 2973
$ws.Range("M113").Value = -7280
$ws.Range("N113").Value = -7313

$ws.Range("H121").Value = 38830.52
$ws.Range("I121").Value = 111611.78
$ws.Range("K121").Value = 334835.34
$ws.Range("M121").Value = -333525.34

$ws.Range("H122").Value = 1711.6923
$ws.Range("I122").Value = 883
$ws.Range("J122").Value = 2080
$ws.Range("K122").Value = 7947
$ws.Range("L122").Value = 18720
$ws.Range("M122").Value = -5497
$ws.Range("N122").Value = -23620

$ws.Range("H135").Value = 601.13336
$ws.Range("I135").Value = 601.13336
$ws.Range("K135").Value = 5410.20024
$ws.Range("M135").Value = -2875.20024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6019.1665
$ws.Range("I2").Value = 7712.3076
$ws.Range("J2").Value = 1617
$ws.Range("K2").Value = 7712.3076
$ws.Range("L2").Value = 1617
$ws.Range("M2").Value = -7599.3076
$ws.Range("N2").Value = -1843

$ws.Range("H9").Value = 656.8570999999999
$ws.Range("J9").Value = 529.6
$ws.Range("L9").Value = 529.6
$ws.Range("N9").Value = -869.6

$ws.Range("H97").Value = 2534
$ws.Range("I97").Value = 2871.7273
$ws.Range("J97").Value = 1914.8334
$ws.Range("K97").Value = 2871.7273
$ws.Range("L97").Value = 1914.8334
$ws.Range("M97").Value = -2375.7273
$ws.Range("N97").Value = -2906.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 247.44444
$ws.Range("I9").Value = 247.44444
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 247.44444
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -23.44443999999999
$ws.Range("N9").ClearContents()

$ws.Range("H13").Value = 6003
$ws.Range("I13").Value = 5006
$ws.Range("J13").Value = 7000
$ws.Range("K13").Value = 5006
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = -4866
$ws.Range("N13").Value = -7280

$ws.Range("H40").Value = 3611.8667
$ws.Range("I40").Value = 889
$ws.Range("J40").Value = 4030.7693
$ws.Range("K40").Value = 889
$ws.Range("L40").Value = 4030.7693
$ws.Range("M40").Value = -753
$ws.Range("N40").Value = -4302.7693

$ws.Range("H63").Value = 335572000
$ws.Range("I63").Value = 28000
$ws.Range("J63").Value = 377515000
$ws.Range("K63").Value = 28000
$ws.Range("L63").Value = 377515000
$ws.Range("M63").Value = -27251
$ws.Range("N63").Value = -377516498

$ws.Range("H66").Value = 335572000
$ws.Range("I66").Value = 28000
$ws.Range("J66").Value = 377515000
$ws.Range("K66").Value = 84000
$ws.Range("L66").Value = 1132545000
$ws.Range("M66").Value = -80256
$ws.Range("N66").Value = -1132552488

$ws.Range("H82").Value = 34484816
$ws.Range("I82").Value = 876.13336
$ws.Range("J82").Value = 71431896
$ws.Range("K82").Value = 876.13336
$ws.Range("L82").Value = 71431896
$ws.Range("M82").Value = -515.13336
$ws.Range("N82").Value = -71432618

$ws.Range("H85").Value = 34484816
$ws.Range("I85").Value = 876.13336
$ws.Range("J85").Value = 71431896
$ws.Range("K85").Value = 876.13336
$ws.Range("L85").Value = 71431896
$ws.Range("M85").Value = 371.86664
$ws.Range("N85").Value = -71434392

$ws.Range("H132").Value = 6499.857
$ws.Range("I132").Value = 7099.8
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 21299.4
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -18769.4
$ws.Range("N132").Value = -20060

$ws.Range("H136").Value = 22519.893
$ws.Range("I136").Value = 4918.3447
$ws.Range("J136").Value = 86325.5
$ws.Range("K136").Value = 14755.0341
$ws.Range("L136").Value = 258976.5
$ws.Range("M136").Value = -12205.0341
$ws.Range("N136").Value = -264076.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 27778224
$ws.Range("I107").Value = 407.15384
$ws.Range("J107").Value = 100000550
$ws.Range("K107").Value = 1221.46152
$ws.Range("L107").Value = 300001650
$ws.Range("M107").Value = 698.5384799999999
$ws.Range("N107").Value = -300005490

$ws.Range("H126").Value = 11907829
$ws.Range("I126").Value = 2806
$ws.Range("J126").Value = 27781194
$ws.Range("K126").Value = 8418
$ws.Range("L126").Value = 83343582
$ws.Range("M126").Value = -5948
$ws.Range("N126").Value = -83348522

